# Change people_widget row (row 3 on "survey" sheet) from a
# select_or_add_multiple widget to a textarea widget.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# type: select_or_add_multiple -> textarea
$ws.Range("C3").Value = "textarea"

# value: "Ruben, Jessica" -> "Ruben<newline>Jessica"
$ws.Range("M3").Value = "`"Ruben`nJessica`""

# M3 should now wrap its (now multi-line) text, matching style used
# elsewhere in the row (style index 7, which has wrapText enabled)
$ws.Range("M3").WrapText = $true

# Update the active selection to reflect the new value cell
$ws.Range("M3").Select()
